$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("levers")

# ---- Read current C, D values (rows 2..21) ----
$c = @{}
$d = @{}
for ($r = 2; $r -le 21; $r++) {
    $c[$r] = $ws.Cells.Item($r, 3).Value2
    $d[$r] = $ws.Cells.Item($r, 4).Value2
}

# ---- Shift C column down by one row (2..20 -> 3..21); old row21 -> row2 ----
for ($r = 21; $r -ge 3; $r--) {
    $ws.Cells.Item($r, 3).Value = $c[$r-1]
}
$ws.Cells.Item(2, 3).Value = $c[21]

# ---- D column: rows 2 and 3 become plain values; rows 4..21 become formulas D{r-1}+0.5 ----
$ws.Cells.Item(2, 4).Value = $d[21]
$ws.Cells.Item(3, 4).Value = $d[2]
for ($r = 4; $r -le 21; $r++) {
    $ws.Cells.Item($r, 4).Formula = "=D" + ($r-1) + "+0.5"
}

# ---- B column: formula "ADV-"&D{r} for every row (relative to its own row) ----
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 2).Formula = '="ADV-"&D' + $r
}

Write-Host "done"
